$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Farmacias")
$ws.Range("I3").Value = "(11) 4961-0338"
$ws.Range("I8").Value = "(11) 4505-1010"
[void]$ws.Range("H19").Select()
